$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the cells we are about to rewrite to Text format so the new values are
# stored as literal strings (matching the inlineStr text cells already used in
# the source file) rather than being auto-coerced into numbers / percentages.
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D47:E47").NumberFormat = "@"

$ws.Range("D2").Value = "287.31"
$ws.Range("E2").Value = "2.11%"
$ws.Range("G2").Value = "16"

$ws.Range("D3").Value = "29.30"
$ws.Range("E3").Value = "3.82%"
$ws.Range("G3").Value = "16"

$ws.Range("D4").Value = "5.094"
$ws.Range("E4").Value = "0.88%"
$ws.Range("G4").Value = "16"

$ws.Range("D5").Value = "0.06997"
$ws.Range("E5").Value = "7.96%"
$ws.Range("G5").Value = "16"

$ws.Range("D6").Value = "7.414"
$ws.Range("E6").Value = "2.38%"
$ws.Range("G6").Value = "16"

$ws.Range("D7").Value = "3.569"
$ws.Range("E7").Value = "5.63%"
$ws.Range("G7").Value = "16"

$ws.Range("D8").Value = "1.429"
$ws.Range("E8").Value = "3.53%"
$ws.Range("G8").Value = "16"

$ws.Range("D9").Value = "0.8993"
$ws.Range("E9").Value = "-3.50%"
$ws.Range("G9").Value = "16"

$ws.Range("E10").Value = "2.57%"
$ws.Range("G10").Value = "16"

$ws.Range("D11").Value = "0.07205"
$ws.Range("E11").Value = "20.49%"
$ws.Range("G11").Value = "16"

$ws.Range("E12").Value = "1.94%"
$ws.Range("G12").Value = "16"

$ws.Range("D13").Value = "0.02929"
$ws.Range("E13").Value = "0.52%"
$ws.Range("G13").Value = "16"

$ws.Range("D14").Value = "0.08998"
$ws.Range("E14").Value = "0.16%"
$ws.Range("G14").Value = "16"

$ws.Range("D15").Value = "0.001614"
$ws.Range("E15").Value = "1.57%"
$ws.Range("G15").Value = "16"

$ws.Range("D16").Value = "0.0006477"
$ws.Range("E16").Value = "1.60%"
$ws.Range("G16").Value = "16"

$ws.Range("D17").Value = "0.006383"
$ws.Range("E17").Value = "4.23%"
$ws.Range("G17").Value = "16"

$ws.Range("D18").Value = "3.470"
$ws.Range("E18").Value = "0.40%"
$ws.Range("G18").Value = "16"

$ws.Range("D19").Value = "2.231"
$ws.Range("E19").Value = "-0.13%"
$ws.Range("G19").Value = "16"

$ws.Range("D20").Value = "0.3233"
$ws.Range("E20").Value = "1.04%"
$ws.Range("G20").Value = "16"

$ws.Range("D21").Value = "0.1325"
$ws.Range("E21").Value = "1.59%"
$ws.Range("G21").Value = "16"

$ws.Range("D22").Value = "4.012"
$ws.Range("E22").Value = "-2.13%"
$ws.Range("G22").Value = "16"

$ws.Range("D23").Value = "0.1557"
$ws.Range("E23").Value = "0.66%"
$ws.Range("G23").Value = "16"

$ws.Range("D24").Value = "0.04533"
$ws.Range("E24").Value = "2.24%"
$ws.Range("G24").Value = "16"

$ws.Range("D25").Value = "0.001208"
$ws.Range("E25").Value = "1.89%"
$ws.Range("G25").Value = "16"

$ws.Range("D26").Value = "0.004379"
$ws.Range("E26").Value = "-0.12%"
$ws.Range("G26").Value = "16"

$ws.Range("D27").Value = "0.0001169"
$ws.Range("E27").Value = "-6.57%"
$ws.Range("G27").Value = "16"

$ws.Range("D28").Value = "0.0001616"
$ws.Range("E28").Value = "-0.24%"
$ws.Range("G28").Value = "16"

$ws.Range("G29").Value = "16"

$ws.Range("G30").Value = "16"

$ws.Range("G31").Value = "16"

$ws.Range("G32").Value = "16"

$ws.Range("G33").Value = "16"

$ws.Range("G34").Value = "16"

$ws.Range("G35").Value = "16"

$ws.Range("G36").Value = "16"

$ws.Range("G37").Value = "16"

$ws.Range("G38").Value = "16"

$ws.Range("G39").Value = "16"

$ws.Range("D40").Value = "0.04270"
$ws.Range("E40").Value = "2.89%"
$ws.Range("G40").Value = "16"

$ws.Range("D41").Value = "0.006807"
$ws.Range("E41").Value = "2.37%"
$ws.Range("G41").Value = "16"

$ws.Range("D42").Value = "0.1250"
$ws.Range("E42").Value = "2.37%"
$ws.Range("G42").Value = "16"

$ws.Range("D43").Value = "0.002107"
$ws.Range("E43").Value = "3.75%"
$ws.Range("G43").Value = "16"

$ws.Range("D44").Value = "0.01179"
$ws.Range("E44").Value = "-2.51%"
$ws.Range("G44").Value = "16"

$ws.Range("D45").Value = "0.00005787"
$ws.Range("E45").Value = "4.89%"
$ws.Range("G45").Value = "16"

$ws.Range("G46").Value = "16"

$ws.Range("D47").Value = "0.01306"
$ws.Range("E47").Value = "0.29%"
$ws.Range("G47").Value = "16"

$ws.Range("G48").Value = "16"

$ws.Range("G49").Value = "16"

$ws.Range("G50").Value = "16"

$ws.Range("G51").Value = "16"
